$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.918.13'
$ws.Range('E2').Value = '  +6.39%  '
$ws.Range('D3').Value = '1.733.82'
$ws.Range('E3').Value = '  +4.65%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9991'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '229.37'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.46%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5438'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.89%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9989'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2776'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +4.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06788'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +6.61%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.72'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.31%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07809'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.720'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.52%  '
$ws.Range('D13').Value = '1.724.28'
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('D14').Value = '1.961.34'
$ws.Range('E14').Value = '  +4.09%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6021'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +6.72%  '
$ws.Range('D16').Value = '0.0₅8412'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '68.73'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +5.00%  '
$ws.Range('D18').Value = '27.835.54'
$ws.Range('E18').Value = '  +6.14%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '215.45'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +12.03%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.848'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.31%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9996'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.93'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +5.25%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.259'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.9996'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.56'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1250'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +4.12%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.453'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.47%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.639'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +9.35%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '16.87'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +5.72%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05611'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.316'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.84%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.712'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.05%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.535'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.640'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.70%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9829'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.91%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.863'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.433'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5951'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.55%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01669'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.53%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.909'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.044.75'
$ws.Range('E41').Value = '  +2.12%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8419'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9985'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '102.78'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('D45').Value = '1.870.24'
$ws.Range('E45').Value = '  +4.18%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '59.91'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈111'
$ws.Range('E47').Value = '  +4.51%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.276'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.005'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.4409'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05299'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.36%  '
